# Apply updated TOTAL_SUBSTATION_LOAD (B), CONTESTABLE_ENERGY (C) and
# ACTUAL_ENERGY (D) figures to the data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @(B, C, D) ; $null means "leave cell as-is / do not create it"
$data = @{
    2  = @(36728, 5799.800999999999, 30928.199)
    3  = @(35269, 5649.136500000001, 29619.8635)
    4  = @(33348, 5554.8325, 27793.1675)
    5  = @(31700, 5481.98, 26218.02)
    6  = @(30762, 5456.8815, 25305.1185)
    7  = @(30953, 5512.395, 25440.605)
    8  = @(31330, 5765.477227722772, 25564.52277227723)
    9  = @(32595, 6519.4465, 26075.5535)
    10 = @(36208, 7901.5895, 28306.4105)
    11 = @(39069, 13052.0005, 26016.9995)
    12 = @(39884, 14886.7985, 24997.2015)
    13 = @(39218, 14810.005, 24407.995)
    14 = @(38683, 14699.3385, 23983.6615)
    15 = @(40878, 15357.636, 25520.364)
    16 = @(41508, 15432.207, 26075.793)
    17 = @(40766, 15338.9075, 25427.0925)
    18 = @(38994, 15792.80217625723, 23201.19782374277)
    19 = @(36246, 15557.50171551809, 20688.49828448191)
    20 = @($null, 15152.49877462994, $null)
    21 = @($null, 13812.90902852661, $null)
    22 = @($null, 12007.84432898735, $null)
    23 = @($null, 9523.143, $null)
    24 = @($null, 6375.7855, $null)
    25 = @($null, 5494.996500000001, $null)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $bVal = $vals[0]
    $cVal = $vals[1]
    $dVal = $vals[2]

    if ($null -ne $bVal) {
        $ws.Cells.Item($row, 2).Value = $bVal
    }
    if ($null -ne $cVal) {
        $ws.Cells.Item($row, 3).Value = $cVal
    }
    if ($null -ne $dVal) {
        $ws.Cells.Item($row, 4).Value = $dVal
    }
}
